# Append the 15.9.2025 working-hours entry as row 7 (Date / Start / End / Start / End)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7: new date label -> becomes a new shared string "15.9.2025"
$ws.Range("A7").Value = "15.9.2025"

# Match the time formatting used by the existing log rows (B2:E2) before
# writing the new time-of-day values, so the new cells share style index 1.
$ws.Range("B7:E7").NumberFormat = $ws.Range("B2:E2").NumberFormat

$ws.Range("B7").Value = 0.8125                 # 19:30
$ws.Range("C7").Value = 0.85416666666666663    # 20:30
$ws.Range("D7").Value = 0.88888888888888884    # 21:20
$ws.Range("E7").Value = 0.010416666666666666   # 00:15 (past midnight)

# Move/save the active selection to F9, matching the edited workbook state.
$ws.Range("F9").Select()
